# Fruta / hortaliza, semanal
#
# Insert a new weekly price record for Mango at row 117, pushing the
# existing records (old rows 117-210) down by one row to 118-211.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 117 (shifts rows 117:210 down to 118:211)
$ws.Rows("117:117").Insert()

# Populate the newly inserted row 117 with the new observation
$ws.Range("A117").Value = 7
$ws.Range("B117").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C117").Value = "Ñuble"
$ws.Range("D117").Value = 45240
$ws.Range("E117").Value = 16
$ws.Range("F117").Value = "Fruta"
$ws.Range("G117").Value = 100108
$ws.Range("H117").Value = "Tropicales y subtropicales"
$ws.Range("I117").Value = 100108002
$ws.Range("J117").Value = "Mango"
$ws.Range("K117").Value = "Sin especificar"
$ws.Range("L117").Value = "Primera"
$ws.Range("M117").Value = 100
$ws.Range("N117").Value = 12000
$ws.Range("O117").Value = 12000
$ws.Range("P117").Value = 12000
$ws.Range("Q117").Value = "$/bandeja 4 kilos"
$ws.Range("R117").Value = "Brasil"
$ws.Range("S117").Value = 3000
$ws.Range("T117").Value = 4
